$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Verify item can be deleted from cart and price is decreased"
$ws.Range("B4").Value = "iphone"

$ws.Range("A5").Value = 'Verify items "No Results" is displayed for invalid product name'
$ws.Range("B5").Value = "qqqqqqqqqqqq"

$ws.Range("A6").Value = "Verify add to cart button"
$ws.Range("B6").Value = "iphone"

$ws.Range("B6").Select()
